{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is executed as: async (context) => { ...this... }\n//\n// Goal (per the OOXML diff / commit \"Update artefato 17 e 19\"):\n//  - Paragraph \"RN-0001:\" loses its old body text (and the _GoBack bookmark\n//    that sat inside it) and gets new body text about informing the client\n//    of the service value before payment.\n//  - Paragraph \"RN-0002:\" gets a new pPr/rPr (bCs) and new body text about\n//    the employee calling the client to perform the service.\n//  - A brand new paragraph \"RN-0003:\" is appended after it, with body text\n//    about paying at the register when buying a product; the _GoBack\n//    bookmark (id 0) is now placed at the end of this new paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the RN-0001 and RN-0002 paragraphs by their (stable) bold label\n// prefix rather than by a hard-coded index, so the script is resilient to\n// unrelated paragraphs being present before them.\nlet rn1Para = null;\nlet rn2Para = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (rn1Para === null && t.indexOf(\"RN-0001:\") === 0) {\n    rn1Para = paragraphs.items[i];\n  } else if (rn2Para === null && t.indexOf(\"RN-0002:\") === 0) {\n    rn2Para = paragraphs.items[i];\n  }\n}\n\nif (!rn1Para || !rn2Para) {\n  throw new Error(\"Could not locate RN-0001/RN-0002 paragraphs\");\n}\n\nconst pkgOpen =\n  '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>';\nconst pkgClose = \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n// --- RN-0001 paragraph: replace its whole range (this also drops the old\n// _GoBack bookmark that used to live inside it). ---\nconst rn1Ooxml =\n  pkgOpen +\n  \"<w:p>\" +\n  \"<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-0001:</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> Para finalizar o pedido, o valor do servi\u00e7o deve ser informado ao cliente para o cliente realizar o pagamento.</w:t></w:r>' +\n  \"</w:p>\" +\n  pkgClose;\n\nrn1Para.getRange().insertOoxml(rn1Ooxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- RN-0002 paragraph: replace its whole range AND append the brand new\n// RN-0003 paragraph right after it (both come from a single OOXML package\n// with two <w:p> elements, which \"Replace\" turns into two real\n// paragraphs). The _GoBack bookmark is re-created at the end of RN-0003.\n// ---\nconst rn2AndRn3Ooxml =\n  pkgOpen +\n  \"<w:p>\" +\n  \"<w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-000</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2</w:t></w:r>\" +\n  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n  \"<w:r><w:rPr><w:bCs/></w:rPr><w:t>Para a realiza\u00e7\u00e3o do servi\u00e7o, o funcion\u00e1rio deve chamar o cliente para</w:t></w:r>\" +\n  '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> a</w:t></w:r>' +\n  '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> realiza</w:t></w:r>' +\n  \"<w:r><w:rPr><w:bCs/></w:rPr><w:t>\u00e7\u00e3o</w:t></w:r>\" +\n  '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  \"<w:r><w:rPr><w:bCs/></w:rPr><w:t>d</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:bCs/></w:rPr><w:t>o servi\u00e7o</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:bCs/></w:rPr><w:t>.</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-000</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3</w:t></w:r>\" +\n  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n  '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=\"preserve\">Ao comprar um produto, </w:t></w:r>' +\n  \"<w:r><w:rPr><w:bCs/></w:rPr><w:t>o cliente</w:t></w:r>\" +\n  '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> deve pagar junto ao caixa.</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  \"</w:p>\" +\n  pkgClose;\n\nrn2Para.getRange().insertOoxml(rn2AndRn3Ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Goal (per the OOXML diff / commit \"Update artefato 17 e 19\"):\n#  - Paragraph \"RN-0001:\" loses its old body text (and the _GoBack bookmark\n#    that sat inside it) and gets new body text about informing the client\n#    of the service value before payment.\n#  - Paragraph \"RN-0002:\" gets a new pPr/rPr (bCs) and new body text about\n#    the employee calling the client to perform the service.\n#  - A brand new paragraph \"RN-0003:\" is appended after it, with body text\n#    about paying at the register when buying a product; the _GoBack\n#    bookmark (id 0) is now placed at the end of this new paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the RN-0001 and RN-0002 paragraphs by their (stable) bold label\n# prefix rather than by a hard-coded index, so the script is resilient to\n# unrelated paragraphs being present before them.\n$rn1Para = $null\n$rn2Para = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($rn1Para -eq $null -and $t.StartsWith(\"RN-0001:\")) {\n        $rn1Para = $p\n    } elseif ($rn2Para -eq $null -and $t.StartsWith(\"RN-0002:\")) {\n        $rn2Para = $p\n    }\n}\n\nif ($rn1Para -eq $null -or $rn2Para -eq $null) {\n    throw \"Could not locate RN-0001/RN-0002 paragraphs\"\n}\n\n$pkgOpen = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# --- RN-0001 paragraph: replace its whole range (this also drops the old\n# _GoBack bookmark that used to live inside it). ---\n$rn1Ooxml = $pkgOpen +\n    '<w:p>' +\n    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-0001:</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> Para finalizar o pedido, o valor do servi\u00e7o deve ser informado ao cliente para o cliente realizar o pagamento.</w:t></w:r>' +\n    '</w:p>' +\n    $pkgClose\n\n$null = $rn1Para.Range.InsertXML($rn1Ooxml)\n\n# --- RN-0002 paragraph: replace its whole range AND append the brand new\n# RN-0003 paragraph right after it (both come from a single OOXML package\n# with two <w:p> elements, which InsertXML turns into two real\n# paragraphs). The _GoBack bookmark is re-created at the end of RN-0003.\n# ---\n$rn2AndRn3Ooxml = $pkgOpen +\n    '<w:p>' +\n    '<w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-000</w:t></w:r>' +\n    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2</w:t></w:r>' +\n    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t>Para a realiza\u00e7\u00e3o do servi\u00e7o, o funcion\u00e1rio deve chamar o cliente para</w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> a</w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> realiza</w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t>\u00e7\u00e3o</w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t>d</w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t>o servi\u00e7o</w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t>.</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p>' +\n    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-000</w:t></w:r>' +\n    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3</w:t></w:r>' +\n    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=\"preserve\">Ao comprar um produto, </w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t>o cliente</w:t></w:r>' +\n    '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> deve pagar junto ao caixa.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>' +\n    $pkgClose\n\n$null = $rn2Para.Range.InsertXML($rn2AndRn3Ooxml)\n"}
